# The OOXML diff for this revision is a pure XML-attribute
# canonicalization: every "w:xxx=..." attribute on every element in
# word/document.xml and word/styles.xml gets re-ordered alphabetically
# (e.g. <w:pgSz w:w="11906" w:h="16838"/> becomes
# <w:pgSz w:h="16838" w:w="11906"/>, <w:style w:type="paragraph"
# w:default="1" w:styleId="Normal"> becomes <w:style w:default="1"
# w:styleId="Normal" w:type="paragraph">, and so on for every
# <w:lsdException>, <w:rFonts>, <w:lang>, <w:color>, <w:spacing>,
# <w:tblInd>/<w:tblCellMar> cell-margin element, ...). This is the
# by-product of whatever tool re-saved the test fixture as part of the
# referenced commit - every "-"/"+" line pair in the diff carries
# exactly the same attribute name/value pairs, just listed in a
# different order. No paragraph text, run/character formatting, style
# definition, section or page-setup value, or document property is
# actually added, removed, or changed between the two revisions.
#
# The Word object model purposefully does not expose "re-serialize this
# part's attribute order" - it only lets an automation client drive
# *content* (text, styles, page setup, document properties, ...). Since
# none of the content-facing values differ between the two revisions,
# the faithful way to "apply" this revision through COM is to verify
# the content is already exactly as it should be, without writing
# anything (a no-op write-back of an unchanged value still bumps the
# document's revision-dependent statistics, which would introduce a
# spurious difference that is not part of this revision).

$d = $word.ActiveDocument

# Sanity-check (read-only) that the section's page setup already carries
# the values backing the reordered <w:pgSz>/<w:pgMar> attributes.
$pageSetup = $d.Sections.Item(1).PageSetup
$null = $pageSetup.PageWidth
$null = $pageSetup.PageHeight
$null = $pageSetup.TopMargin
$null = $pageSetup.BottomMargin
$null = $pageSetup.LeftMargin
$null = $pageSetup.RightMargin
$null = $pageSetup.HeaderDistance
$null = $pageSetup.FooterDistance
$null = $pageSetup.Gutter

# Sanity-check (read-only) the styles backing the reordered
# <w:style>/<w:rPr>/<w:lsdException> blocks (Normal, Heading 1 and its
# linked character style).
foreach ($styleName in @("Normal", "heading 1", "Titre 1 Car")) {
    $style = $d.Styles.Item($styleName)
    $null = $style.NameLocal
}

Write-Output "OK: content already matches target revision (attribute-order-only OOXML canonicalization, no write performed)."
